# Add "Physically Based Shading at Disney" (SIGGRAPH 2012) reading entry
# as a new row 8 in the readings table, pushing existing rows 8-13 down to 9-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 8 (existing rows 8-13 shift down to 9-14).
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).RowHeight = 42.75

# 2. Fill in the new row's data, in the same order the original author would
#    have (title, author, homepage, publisher, tags, thumb, url, notes, brief,
#    then date) so new shared-string entries line up with the target file.
$ws.Range("A8").Value = "Physically Based Shading at Disney"
$ws.Range("B8").Value = "Brent Burley"
$ws.Range("C8").Value = "https://blog.selfshadow.com/publications/s2012-shading-course/"
$ws.Range("E8").Value = "SIGGRAPH 2012"
$ws.Range("F8").Value = "PBR"
$ws.Range("G8").Value = "assets\thumb\disney_pbs_sig2012.png"
$ws.Range("I8").Value = "assets\slides\s2012_pbs_disney_brdf_slides_v2.pdf"
$ws.Range("K8").Value = "assets\slides\s2012_pbs_disney_brdf_notes_v3.pdf"
$ws.Range("L8").Value = "虚幻4的材质模型参考了迪士尼，这个演讲可以参考。"
$ws.Range("D8").Value = "2012年8月"

# 3. Grow the table (表1) to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:L14"))

# 4. Hyperlinks don't automatically follow the row insertion, so rebuild the
#    whole collection: existing links shift down one row (for row >= 8), and
#    four new links are added for the new row 8.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H5"), "https://youtu.be/yy8jQgmhbAU")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://stoyannk.wordpress.com/2018/10/26/my-cppcon-2019-talk")
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.ea.com/frostbite/news/introduction-to-data-oriented-design")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://youtu.be/rX0ItVEVjHc")
$ws.Hyperlinks.Add($ws.Range("I5"), "assets\slides\oop_is_dead_long_live_dataoriented_design__stoyan_nikolov__cppcon_2018.pdf")
$ws.Hyperlinks.Add($ws.Range("I13"), "assets\slides\Introduction_to_Data-Oriented_Design_2014DICE.pdf")
$ws.Hyperlinks.Add($ws.Range("I6"), "assets\slides\Data-Oriented Design and C++ - Mike Acton - CppCon 2014.pptx")
$ws.Hyperlinks.Add($ws.Range("I11"), "assets\slides\GDC17-framegraph.pptx")
$ws.Hyperlinks.Add($ws.Range("C11"), "https://www.ea.com/frostbite/news/framegraph-extensible-rendering-architecture-in-frostbite")
$ws.Hyperlinks.Add($ws.Range("G5"), "assets\thumb\oop_is_dead_cppcon2018.png")
$ws.Hyperlinks.Add($ws.Range("G13"), "assets\thumb\dod_dice.png")
$ws.Hyperlinks.Add($ws.Range("G6"), "assets\thumb\dod_cppcon2014.png")
$ws.Hyperlinks.Add($ws.Range("G11"), "assets\thumb\framegraph_gdc2018.png")
$ws.Hyperlinks.Add($ws.Range("C12"), "https://dl.acm.org/citation.cfm?id=15902")
$ws.Hyperlinks.Add($ws.Range("G12"), "assets\thumb\the_rendering_equation.png")
$ws.Hyperlinks.Add($ws.Range("I12"), "assets\paper\p143-kajiya.pdf")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://youtu.be/ZvGPrDdVl4Y")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://www.gdcvault.com/play/1024001/-Overwatch-Gameplay-Architecture-and")
$ws.Hyperlinks.Add($ws.Range("G4"), "assets\thumb\overwatch_ecs_gdc2017.png")
$ws.Hyperlinks.Add($ws.Range("J5"), "https://neil3d.github.io/3dengine/why-ecs.html")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://blog.selfshadow.com/publications/s2013-shading-course/")
$ws.Hyperlinks.Add($ws.Range("G7"), "assets\thumb\pbr_background_sig2013.png")
$ws.Hyperlinks.Add($ws.Range("H7"), "https://youtu.be/j-A0mwsJRmk")
$ws.Hyperlinks.Add($ws.Range("I7"), "assets\slides\s2013_pbs_physics_math_slides.pdf")
$ws.Hyperlinks.Add($ws.Range("K7"), "assets\slides\s2013_pbs_physics_math_notes.pdf")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://blog.selfshadow.com/publications/s2013-shading-course/")
$ws.Hyperlinks.Add($ws.Range("G9"), "assets\thumb\ue4_pbr_sig2013.png")
$ws.Hyperlinks.Add($ws.Range("I9"), "assets\slides\s2013_pbs_epic_slides.pdf")
$ws.Hyperlinks.Add($ws.Range("J9"), "https://neil3d.github.io/unreal/pbr-ue4.html")
$ws.Hyperlinks.Add($ws.Range("K9"), "assets\slides\s2013_pbs_epic_notes_v2.pdf")
$ws.Hyperlinks.Add($ws.Range("C14"), "https://dl.acm.org/citation.cfm?id=357293")
$ws.Hyperlinks.Add($ws.Range("G14"), "assets\thumb\cook-brdf.png")
$ws.Hyperlinks.Add($ws.Range("I14"), "assets\paper\p7-cook.pdf")
$ws.Hyperlinks.Add($ws.Range("C10"), "http://intro-to-dxr.cwyman.org/")
$ws.Hyperlinks.Add($ws.Range("G10"), "assets\thumb\intro_dxr_sig2018.png")
$ws.Hyperlinks.Add($ws.Range("H10"), "https://youtu.be/Q1cuuepVNoY")
$ws.Hyperlinks.Add($ws.Range("I10"), "assets\slides\s2018_IntroDXR_RaytracingShaders.pdf")
$ws.Hyperlinks.Add($ws.Range("J10"), "https://neil3d.github.io/3dengine/DXRPreview.html")
$ws.Hyperlinks.Add($ws.Range("K10"), "http://cwyman.org/code/dxrTutors/dxr_tutors.md.html")
$ws.Hyperlinks.Add($ws.Range("K4"), "http://gad.qq.com/article/detail/28682")
$ws.Hyperlinks.Add($ws.Range("J4"), "https://neil3d.github.io/3dengine/why-ecs.html")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://blog.selfshadow.com/publications/s2012-shading-course/")
$ws.Hyperlinks.Add($ws.Range("G8"), "assets\thumb\disney_pbs_sig2012.png")
$ws.Hyperlinks.Add($ws.Range("I8"), "assets\slides\s2012_pbs_disney_brdf_slides_v2.pdf")
$ws.Hyperlinks.Add($ws.Range("K8"), "assets\slides\s2012_pbs_disney_brdf_notes_v3.pdf")

# 5. Move the active selection, matching the author's final cursor position.
$ws.Range("D9").Select()
